# "data discussion slide order"
# Swap the order of the two "DISCUSSION QUESTION" slides:
#   - "How could record linkage choices and methods impact bias and fairness?"
#   - "What are some experiences and challenges you've encountered with acquiring data ...?"
# The record-linkage slide (originally at position 10) moves after the
# acquiring-data-experiences slide (originally at position 11).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(10)
$s.MoveTo(11)
